# Daily attendance processing - 2025-10-22 08:54:12
# Normalizes the "Recorded By" column (G) order for each attendance row:
# when a cell holds a comma-separated list of recorders, the first two
# entries are swapped - except for the already-canonical
# "System, backup@backdoor.com" pairing, which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value()

    if ($v -ne $null) {
        if ($v -ne "System, backup@backdoor.com") {
            $parts = $v.Split(",")

            if ($parts.Count -ge 2) {
                $trimmed = @()
                foreach ($p in $parts) {
                    $trimmed += $p.Trim()
                }

                $first = $trimmed[0]
                $second = $trimmed[1]
                $trimmed[0] = $second
                $trimmed[1] = $first

                $newVal = $trimmed[0]
                for ($i = 1; $i -lt $trimmed.Count; $i++) {
                    $newVal = $newVal + ", " + $trimmed[$i]
                }

                $cell.Value = $newVal
            }
        }
    }
}
